# Insert a new weekly data row at row 64, pushing existing rows 64-97 down to 65-98.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before the current row 64 (shifts rows 64-97 down to 65-98).
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new data point.
$ws.Cells.Item(64, 1).Value = 11
$ws.Cells.Item(64, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(64, 3).Value = "Bíobío"
$ws.Cells.Item(64, 4).Value = 44603
$ws.Cells.Item(64, 4).Style = $ws.Cells.Item(65, 4).Style
$ws.Cells.Item(64, 4).NumberFormat = $ws.Cells.Item(65, 4).NumberFormat
$ws.Cells.Item(64, 5).Value = 8
$ws.Cells.Item(64, 6).Value = 100112032
$ws.Cells.Item(64, 7).Value = "Zapallo italiano"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 270
$ws.Cells.Item(64, 11).Value = 7500
$ws.Cells.Item(64, 12).Value = 8000
$ws.Cells.Item(64, 13).Value = 7722
$ws.Cells.Item(64, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(64, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(64, 16).Value = 129
$ws.Cells.Item(64, 17).Value = 60
$ws.Cells.Item(64, 18).Value = "Hortaliza"
